# "Type Providers.pptx" — oredev final changes
#
# Slide 2 ("Why F#?") holds a SmartArt diagram (Content Placeholder 5).
# Two of its leaf nodes swap their text:
#   "Parallelization"            -> "Interoperability with .NET"
#   "Interoperability with .NET" -> "Parallelization"
#
# Editing through the SmartArt object model (Shape.SmartArt.AllNodes)
# updates both the cached diagram drawing (dsp:txBody) and the diagram
# data model (dgm:pt) parts together, the same way PowerPoint's UI does.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$sa = $sh.SmartArt

for ($i = 1; $i -le $sa.AllNodes.Count; $i++) {
    $node = $sa.AllNodes.Item($i)
    $tr = $node.TextFrame2.TextRange
    if ($tr.Text -eq "Parallelization") {
        $tr.Text = "Interoperability with .NET"
    }
    elseif ($tr.Text -eq "Interoperability with .NET") {
        $tr.Text = "Parallelization"
    }
}
